$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 7 (ComCast Business): mark as paid/used ("x") with a new bill amount.
$ws.Range("B7").Value = "x"
$ws.Range("D7").Value = 504.2

# Row 15 (Lakeshore Recyling): clear out the old, now-stale invoice entry
# (checkbox, invoice number, bill amount) - keeping only the vendor name.
$ws.Range("B15").ClearContents()
$ws.Range("C15").ClearContents()
$ws.Range("D15").ClearContents()

# Update the tracking/reference number in I1.
$ws.Range("I1").Value = 33404

# Move the active selection to I1.
[void]$ws.Range("I1").Select()
